# Update countries & provincias Spain
# Applies the COVID-19 data refresh (8 Oct 2020, 14:02 -> 15:19) to the
# "Pais" sheet: updated totals for several countries, plus the swapped
# Islas Malvinas / Montserrat rows, and the refreshed timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 7780586
$ws.Cells.Item(4, 3).Value = 4362
$ws.Cells.Item(4, 4).Value = 4984528
$ws.Cells.Item(4, 5).Value = 2579169
$ws.Cells.Item(4, 7).Value = 105
$ws.Cells.Item(4, 8).Value = 216889

# Row 18 - Irak
$ws.Cells.Item(18, 2).Value = 394566
$ws.Cells.Item(18, 3).Value = 3522
$ws.Cells.Item(18, 4).Value = 323815
$ws.Cells.Item(18, 5).Value = 61068
$ws.Cells.Item(18, 7).Value = 79
$ws.Cells.Item(18, 8).Value = 9683

# Row 20 - Arabia Saudita
$ws.Cells.Item(20, 2).Value = 338132
$ws.Cells.Item(20, 3).Value = 421
$ws.Cells.Item(20, 4).Value = 323769
$ws.Cells.Item(20, 5).Value = 9391
$ws.Cells.Item(20, 7).Value = 25
$ws.Cells.Item(20, 8).Value = 4972

# Row 26 - Alemania (D unchanged, H unchanged)
$ws.Cells.Item(26, 2).Value = 311331
$ws.Cells.Item(26, 3).Value = 218
$ws.Cells.Item(26, 5).Value = 33979

# Row 30 - Paises Bajos (D unchanged, E unchanged)
$ws.Cells.Item(30, 2).Value = 155810
$ws.Cells.Item(30, 3).Value = 5822
$ws.Cells.Item(30, 7).Value = 13
$ws.Cells.Item(30, 8).Value = 6531

# Row 36 - Catar (G unchanged, H unchanged)
$ws.Cells.Item(36, 2).Value = 127394
$ws.Cells.Item(36, 3).Value = 213
$ws.Cells.Item(36, 4).Value = 124327
$ws.Cells.Item(36, 5).Value = 2849

# Row 40 - Kuwait
$ws.Cells.Item(40, 2).Value = 109441
$ws.Cells.Item(40, 3).Value = 698
$ws.Cells.Item(40, 4).Value = 101314
$ws.Cells.Item(40, 5).Value = 7485
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = 642

# Row 46 - Suecia (C,D,E,F,H unchanged)
$ws.Cells.Item(46, 2).Value = 97532
$ws.Cells.Item(46, 7).Value = 5

# Row 70 - Estado de Palestina
$ws.Cells.Item(70, 2).Value = 43256
$ws.Cells.Item(70, 3).Value = 416
$ws.Cells.Item(70, 4).Value = 36584
$ws.Cells.Item(70, 5).Value = 6313
$ws.Cells.Item(70, 7).Value = 4
$ws.Cells.Item(70, 8).Value = 359

# Row 76 - Serbia (D unchanged)
$ws.Cells.Item(76, 2).Value = 34344
$ws.Cells.Item(76, 3).Value = 151
$ws.Cells.Item(76, 5).Value = 2048
$ws.Cells.Item(76, 7).Value = 2
$ws.Cells.Item(76, 8).Value = 760

# Row 78 - Dinamarca
$ws.Cells.Item(78, 2).Value = 31156
$ws.Cells.Item(78, 3).Value = 446
$ws.Cells.Item(78, 4).Value = 24899
$ws.Cells.Item(78, 5).Value = 5592
$ws.Cells.Item(78, 7).Value = 2
$ws.Cells.Item(78, 8).Value = 665

# Row 79 - El Salvador (F,G,H unchanged)
$ws.Cells.Item(79, 2).Value = 29842
$ws.Cells.Item(79, 3).Value = 105
$ws.Cells.Item(79, 4).Value = 24770
$ws.Cells.Item(79, 5).Value = 4195

# Row 80 - Bosnia y Herzegovina
$ws.Cells.Item(80, 2).Value = 29528
$ws.Cells.Item(80, 3).Value = 453
$ws.Cells.Item(80, 4).Value = 22939
$ws.Cells.Item(80, 5).Value = 5676
$ws.Cells.Item(80, 7).Value = 5
$ws.Cells.Item(80, 8).Value = 913

# Row 95 - Senegal (F,G,H unchanged)
$ws.Cells.Item(95, 2).Value = 15190
$ws.Cells.Item(95, 3).Value = 16
$ws.Cells.Item(95, 4).Value = 13068
$ws.Cells.Item(95, 5).Value = 1809

# Row 136 - Sri Lanka (only D, E change)
$ws.Cells.Item(136, 4).Value = 3278
$ws.Cells.Item(136, 5).Value = 1168

# Rows 215/216 - Islas Malvinas and Montserrat swap their stats
# (row 215 keeps its position but now carries Montserrat's figures and
# label, and vice-versa for row 216), matching the source diff.
$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1

$ws.Cells.Item(216, 1).Value = "Islas Malvinas"
$ws.Cells.Item(216, 4).Value = 13
$ws.Cells.Item(216, 8).Value = 0

# Footer timestamp refresh
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 15:19"
